# "small fix in all templates"
# Rename the first worksheet and switch the active tab/selection.

$wb = $excel.ActiveWorkbook

# Rename "SampleSubmission" -> "ContainerMove"
$wsMove = $wb.Worksheets.Item("SampleSubmission")
$wsMove.Name = "ContainerMove"

# Make the renamed sheet the active/selected tab (was "Index" before),
# and move its selection to B8 (was D32).
$wsMove.Activate()
$wsMove.Range("B8").Select()
